$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "advis codes" document rows (6 & 7) to point at the new
# 1.0.1 versions, replacing the old local ("t:/...") file paths with the
# new Teams share links.
$ws.Range("B6").Value = "Oversigt_adviskoder_HL7_FHIR_1.0.1.doc"
$ws.Range("D6").Value = "https://teams.microsoft.com/l/file/731C0EC3-5C21-45F2-B54A-8F8D5CB1D510?tenantId=4898802a-6a02-437e-beb6-47cb8db2f43d&fileType=docx&objectUrl=https%3A%2F%2Fmedcomtest.sharepoint.com%2Fsites%2FFHIR-ModerniseringafMedComStandarder%2FDelte%20dokumenter%2FAdvis%20om%20sygehusophold%2FOversigt_adviskoder_HL7_FHIR_1.0.1.docx&baseUrl=https%3A%2F%2Fmedcomtest.sharepoint.com%2Fsites%2FFHIR-ModerniseringafMedComStandarder&serviceName=teams&threadId=19:9f589d52d6f34e96ba534075c28d0abc@thread.skype&groupId=6beb4dc7-667a-48c8-9716-18976d8b8329"

$ws.Range("B7").Value = "Overview_advis codes_HL7_FHIR_1.0.1.docx"
$ws.Range("D7").Value = "https://teams.microsoft.com/l/file/B5D7D898-5663-45C6-9D91-D8A602BBA5FA?tenantId=4898802a-6a02-437e-beb6-47cb8db2f43d&fileType=docx&objectUrl=https%3A%2F%2Fmedcomtest.sharepoint.com%2Fsites%2FFHIR-ModerniseringafMedComStandarder%2FDelte%20dokumenter%2FAdvis%20om%20sygehusophold%2FOverview_advis%20codes_HL7_FHIR_1.0.1.docx&baseUrl=https%3A%2F%2Fmedcomtest.sharepoint.com%2Fsites%2FFHIR-ModerniseringafMedComStandarder&serviceName=teams&threadId=19:9f589d52d6f34e96ba534075c28d0abc@thread.skype&groupId=6beb4dc7-667a-48c8-9716-18976d8b8329"

# Move the active selection to A7, matching the saved view state.
$ws.Range("A7").Select()
